$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-15 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-16 Thursday", 2) | Out-Null
$d.Content.Find.Execute("69-20=", $true, $false, $false, $false, $false, $true, 1, $false, "26+0=", 2) | Out-Null
$d.Content.Find.Execute("26+35=", $true, $false, $false, $false, $false, $true, 1, $false, "23+52=", 2) | Out-Null
$d.Content.Find.Execute("9+30=", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=", 2) | Out-Null
$d.Content.Find.Execute("93-3=", $true, $false, $false, $false, $false, $true, 1, $false, "13+28=", 2) | Out-Null
$d.Content.Find.Execute("54-38=", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=", 2) | Out-Null
$d.Content.Find.Execute("59+9=", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=", 2) | Out-Null
$d.Content.Find.Execute("82-52=", $true, $false, $false, $false, $false, $true, 1, $false, "18+10=", 2) | Out-Null
$d.Content.Find.Execute("91-48=", $true, $false, $false, $false, $false, $true, 1, $false, "89-49=", 2) | Out-Null
$d.Content.Find.Execute("26+14=", $true, $false, $false, $false, $false, $true, 1, $false, "44-1=", 2) | Out-Null
$d.Content.Find.Execute("16+31=", $true, $false, $false, $false, $false, $true, 1, $false, "69-8=", 2) | Out-Null
$d.Content.Find.Execute("19+9=", $true, $false, $false, $false, $false, $true, 1, $false, "38+17=", 2) | Out-Null
$d.Content.Find.Execute("95-58=", $true, $false, $false, $false, $false, $true, 1, $false, "28+24=", 2) | Out-Null
$d.Content.Find.Execute("31-14=", $true, $false, $false, $false, $false, $true, 1, $false, "71+26=", 2) | Out-Null
$d.Content.Find.Execute("42+44=", $true, $false, $false, $false, $false, $true, 1, $false, "17+31=", 2) | Out-Null
$d.Content.Find.Execute("65+13=", $true, $false, $false, $false, $false, $true, 1, $false, "33+30=", 2) | Out-Null
$d.Content.Find.Execute("99-57=", $true, $false, $false, $false, $false, $true, 1, $false, "28+5=", 2) | Out-Null
$d.Content.Find.Execute("14-9=", $true, $false, $false, $false, $false, $true, 1, $false, "73-3=", 2) | Out-Null
$d.Content.Find.Execute("28+6=", $true, $false, $false, $false, $false, $true, 1, $false, "49+3=", 2) | Out-Null
$d.Content.Find.Execute("71-14=", $true, $false, $false, $false, $false, $true, 1, $false, "52+22=", 2) | Out-Null
$d.Content.Find.Execute("97-79=", $true, $false, $false, $false, $false, $true, 1, $false, "64+35=", 2) | Out-Null
$d.Content.Find.Execute("30+30=", $true, $false, $false, $false, $false, $true, 1, $false, "76+4=", 2) | Out-Null
$d.Content.Find.Execute("55+41=", $true, $false, $false, $false, $false, $true, 1, $false, "17+26=", 2) | Out-Null
$d.Content.Find.Execute("39+12=", $true, $false, $false, $false, $false, $true, 1, $false, "83-32=", 2) | Out-Null
$d.Content.Find.Execute("65+29=", $true, $false, $false, $false, $false, $true, 1, $false, "27-2=", 2) | Out-Null
$d.Content.Find.Execute("80-65=", $true, $false, $false, $false, $false, $true, 1, $false, "66+17=", 2) | Out-Null
$d.Content.Find.Execute("55-36=", $true, $false, $false, $false, $false, $true, 1, $false, "18+53=", 2) | Out-Null
$d.Content.Find.Execute("7+67=", $true, $false, $false, $false, $false, $true, 1, $false, "22+26=", 2) | Out-Null
$d.Content.Find.Execute("0+38=", $true, $false, $false, $false, $false, $true, 1, $false, "58+33=", 2) | Out-Null
$d.Content.Find.Execute("60-6=", $true, $false, $false, $false, $false, $true, 1, $false, "16-9=", 2) | Out-Null
$d.Content.Find.Execute("46+39=", $true, $false, $false, $false, $false, $true, 1, $false, "15+6=", 2) | Out-Null
$d.Content.Find.Execute("91+2=", $true, $false, $false, $false, $false, $true, 1, $false, "60-13=", 2) | Out-Null
$d.Content.Find.Execute("61-12=", $true, $false, $false, $false, $false, $true, 1, $false, "11+20=", 2) | Out-Null
$d.Content.Find.Execute("52+8=", $true, $false, $false, $false, $false, $true, 1, $false, "90-68=", 2) | Out-Null
$d.Content.Find.Execute("18+48=", $true, $false, $false, $false, $false, $true, 1, $false, "71-60=", 2) | Out-Null
$d.Content.Find.Execute("31+37=", $true, $false, $false, $false, $false, $true, 1, $false, "95-83=", 2) | Out-Null
$d.Content.Find.Execute("64+25=", $true, $false, $false, $false, $false, $true, 1, $false, "96-55=", 2) | Out-Null
$d.Content.Find.Execute("92-53=", $true, $false, $false, $false, $false, $true, 1, $false, "9+17=", 2) | Out-Null
$d.Content.Find.Execute("64-9=", $true, $false, $false, $false, $false, $true, 1, $false, "21-9=", 2) | Out-Null
$d.Content.Find.Execute("34+2=", $true, $false, $false, $false, $false, $true, 1, $false, "28+25=", 2) | Out-Null
$d.Content.Find.Execute("98-61=", $true, $false, $false, $false, $false, $true, 1, $false, "34+31=", 2) | Out-Null
$d.Content.Find.Execute("20+10=", $true, $false, $false, $false, $false, $true, 1, $false, "88+7=", 2) | Out-Null
$d.Content.Find.Execute("63-26=", $true, $false, $false, $false, $false, $true, 1, $false, "41+52=", 2) | Out-Null
$d.Content.Find.Execute("67-35=", $true, $false, $false, $false, $false, $true, 1, $false, "27+36=", 2) | Out-Null
$d.Content.Find.Execute("11-4=", $true, $false, $false, $false, $false, $true, 1, $false, "62+16=", 2) | Out-Null
$d.Content.Find.Execute("95-51=", $true, $false, $false, $false, $false, $true, 1, $false, "21+32=", 2) | Out-Null
$d.Content.Find.Execute("36-12=", $true, $false, $false, $false, $false, $true, 1, $false, "70-42=", 2) | Out-Null
$d.Content.Find.Execute("11+38=", $true, $false, $false, $false, $false, $true, 1, $false, "78+5=", 2) | Out-Null
$d.Content.Find.Execute("84-16=", $true, $false, $false, $false, $false, $true, 1, $false, "87-36=", 2) | Out-Null
$d.Content.Find.Execute("52-46=", $true, $false, $false, $false, $false, $true, 1, $false, "66-27=", 2) | Out-Null
$d.Content.Find.Execute("73-19=", $true, $false, $false, $false, $false, $true, 1, $false, "25+60=", 2) | Out-Null
$d.Content.Find.Execute("51-28=", $true, $false, $false, $false, $false, $true, 1, $false, "59+8=", 2) | Out-Null
$d.Content.Find.Execute("75-11=", $true, $false, $false, $false, $false, $true, 1, $false, "19+24=", 2) | Out-Null
$d.Content.Find.Execute("20+8=", $true, $false, $false, $false, $false, $true, 1, $false, "24+16=", 2) | Out-Null
$d.Content.Find.Execute("75-2=", $true, $false, $false, $false, $false, $true, 1, $false, "65-59=", 2) | Out-Null
$d.Content.Find.Execute("42-9=", $true, $false, $false, $false, $false, $true, 1, $false, "45+48=", 2) | Out-Null
$d.Content.Find.Execute("19+34=", $true, $false, $false, $false, $false, $true, 1, $false, "74-9=", 2) | Out-Null
$d.Content.Find.Execute("14+41=", $true, $false, $false, $false, $false, $true, 1, $false, "65+3=", 2) | Out-Null
$d.Content.Find.Execute("50-26=", $true, $false, $false, $false, $false, $true, 1, $false, "41+2=", 2) | Out-Null
$d.Content.Find.Execute("26+33=", $true, $false, $false, $false, $false, $true, 1, $false, "79-72=", 2) | Out-Null
$d.Content.Find.Execute("65+34=", $true, $false, $false, $false, $false, $true, 1, $false, "2+89=", 2) | Out-Null
$d.Content.Find.Execute("91-2=", $true, $false, $false, $false, $false, $true, 1, $false, "16+25=", 2) | Out-Null
$d.Content.Find.Execute("28+68=", $true, $false, $false, $false, $false, $true, 1, $false, "68-60=", 2) | Out-Null
$d.Content.Find.Execute("98-7=", $true, $false, $false, $false, $false, $true, 1, $false, "90-4=", 2) | Out-Null
$d.Content.Find.Execute("21+48=", $true, $false, $false, $false, $false, $true, 1, $false, "85-57=", 2) | Out-Null
$d.Content.Find.Execute("56+23=", $true, $false, $false, $false, $false, $true, 1, $false, "95-7=", 2) | Out-Null
$d.Content.Find.Execute("32+42=", $true, $false, $false, $false, $false, $true, 1, $false, "79-25=", 2) | Out-Null
$d.Content.Find.Execute("81-53=", $true, $false, $false, $false, $false, $true, 1, $false, "8+64=", 2) | Out-Null
$d.Content.Find.Execute("23+22=", $true, $false, $false, $false, $false, $true, 1, $false, "9+73=", 2) | Out-Null
$d.Content.Find.Execute("34+16=", $true, $false, $false, $false, $false, $true, 1, $false, "48-44=", 2) | Out-Null
$d.Content.Find.Execute("16+56=", $true, $false, $false, $false, $false, $true, 1, $false, "9+54=", 2) | Out-Null
$d.Content.Find.Execute("3+88=", $true, $false, $false, $false, $false, $true, 1, $false, "48+26=", 2) | Out-Null
$d.Content.Find.Execute("39-5=", $true, $false, $false, $false, $false, $true, 1, $false, "49-31=", 2) | Out-Null
$d.Content.Find.Execute("9+9=", $true, $false, $false, $false, $false, $true, 1, $false, "87-35=", 2) | Out-Null
$d.Content.Find.Execute("37+2=", $true, $false, $false, $false, $false, $true, 1, $false, "31+24=", 2) | Out-Null
$d.Content.Find.Execute("25+55=", $true, $false, $false, $false, $false, $true, 1, $false, "77+15=", 2) | Out-Null
$d.Content.Find.Execute("25+30=", $true, $false, $false, $false, $false, $true, 1, $false, "95-66=", 2) | Out-Null
$d.Content.Find.Execute("56+26=", $true, $false, $false, $false, $false, $true, 1, $false, "63+31=", 2) | Out-Null
$d.Content.Find.Execute("78-17=", $true, $false, $false, $false, $false, $true, 1, $false, "84-44=", 2) | Out-Null
$d.Content.Find.Execute("83-12=", $true, $false, $false, $false, $false, $true, 1, $false, "80-6=", 2) | Out-Null
$d.Content.Find.Execute("22+32=", $true, $false, $false, $false, $false, $true, 1, $false, "87-21=", 2) | Out-Null
$d.Content.Find.Execute("69+23=", $true, $false, $false, $false, $false, $true, 1, $false, "54-33=", 2) | Out-Null
$d.Content.Find.Execute("23-7=", $true, $false, $false, $false, $false, $true, 1, $false, "80-54=", 2) | Out-Null
$d.Content.Find.Execute("60-30=", $true, $false, $false, $false, $false, $true, 1, $false, "4+35=", 2) | Out-Null
$d.Content.Find.Execute("18+74=", $true, $false, $false, $false, $false, $true, 1, $false, "57+33=", 2) | Out-Null
$d.Content.Find.Execute("92-51=", $true, $false, $false, $false, $false, $true, 1, $false, "89-86=", 2) | Out-Null
$d.Content.Find.Execute("72-48=", $true, $false, $false, $false, $false, $true, 1, $false, "16+53=", 2) | Out-Null
$d.Content.Find.Execute("25+28=", $true, $false, $false, $false, $false, $true, 1, $false, "62-54=", 2) | Out-Null
$d.Content.Find.Execute("91-10=", $true, $false, $false, $false, $false, $true, 1, $false, "30+7=", 2) | Out-Null
$d.Content.Find.Execute("67+13=", $true, $false, $false, $false, $false, $true, 1, $false, "81-35=", 2) | Out-Null
$d.Content.Find.Execute("86-33=", $true, $false, $false, $false, $false, $true, 1, $false, "77-41=", 2) | Out-Null
$d.Content.Find.Execute("12+49=", $true, $false, $false, $false, $false, $true, 1, $false, "50+25=", 2) | Out-Null
$d.Content.Find.Execute("2+10=", $true, $false, $false, $false, $false, $true, 1, $false, "1+5=", 2) | Out-Null
$d.Content.Find.Execute("8+44=", $true, $false, $false, $false, $false, $true, 1, $false, "41+39=", 2) | Out-Null
$d.Content.Find.Execute("44+24=", $true, $false, $false, $false, $false, $true, 1, $false, "87-18=", 2) | Out-Null
$d.Content.Find.Execute("10-4=", $true, $false, $false, $false, $false, $true, 1, $false, "10+42=", 2) | Out-Null
$d.Content.Find.Execute("46+53=", $true, $false, $false, $false, $false, $true, 1, $false, "20+73=", 2) | Out-Null
$d.Content.Find.Execute("40+42=", $true, $false, $false, $false, $false, $true, 1, $false, "3+28=", 2) | Out-Null
$d.Content.Find.Execute("33-13=", $true, $false, $false, $false, $false, $true, 1, $false, "90+6=", 2) | Out-Null
$d.Content.Find.Execute("16+55=", $true, $false, $false, $false, $false, $true, 1, $false, "96-79=", 2) | Out-Null
$d.Content.Find.Execute("37-3=", $true, $false, $false, $false, $false, $true, 1, $false, "86-1=", 2) | Out-Null
